$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New glossary terms appended to the end of the (single) column A list.
$values = @(
    "hyper parameter",
    "Stochastic Weight Averaging",
    "Prompt selection",
    "classifier",
    "Frozen",
    "Class feature (CLS)"
)

$startRow = 58
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Mirror the author's final view/selection state: last cell selected & scrolled into view.
$lastRow = $startRow + $values.Length - 1
$lastCell = $ws.Cells.Item($lastRow, 1)
$lastCell.Select()
$excel.ActiveWindow.ScrollRow = 56
$excel.ActiveWindow.ScrollColumn = 1
